$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "'2026-02-07"
$ws.Range("B38").Value = "'3740240"
$ws.Range("C38").Value = "'20"
$ws.Range("D38").Value = "'1"
